# Actualización automática desde tarea programada
# Corrects the timestamp on row 4 (floating point refinement) and appends
# a new data row (row 5) captured by the scheduled DropControl task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: refine the stored timestamp value ---
$ws.Cells.Item(4, 1).Value = 45873.41691729167

# --- Row 5: new sensor reading appended by the scheduled task ---
$ws.Cells.Item(5, 1).Value = 45873.50022492938
$ws.Cells.Item(5, 1).NumberFormat = $ws.Cells.Item(4, 1).NumberFormat

$ws.Cells.Item(5, 2).Value = 2025
$ws.Cells.Item(5, 3).Value = 15
$ws.Cells.Item(5, 4).Value = 19.11
$ws.Cells.Item(5, 5).Value = 77.90000000000001
$ws.Cells.Item(5, 6).Value = 631.05
$ws.Cells.Item(5, 7).Value = 13.6
$ws.Cells.Item(5, 8).Value = "ESE"
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = "12:00:19"
